$d = $word.ActiveDocument

# Pull the full package OOXML for the document content. Working against this
# string directly lets us control exactly how each <w:t> element is produced
# (in particular, preserving the xml:space="preserve" attribute that the
# target revision keeps on every changed run).
$full = $d.Content.WordOpenXML

function Replace-RunText {
    param(
        [string]$Xml,
        [string]$OldText,
        [string]$NewText
    )

    $oldEsc = [regex]::Escape($OldText)
    $pattern = '<w:t(?: xml:space="preserve")?>' + $oldEsc + '</w:t>'
    $regex = New-Object System.Text.RegularExpressions.Regex($pattern)

    $matchCount = $regex.Matches($Xml).Count
    if ($matchCount -ne 1) {
        throw "Expected exactly one match for old text, got $($matchCount): $OldText"
    }

    $newLiteral = '<w:t xml:space="preserve">' + $NewText + '</w:t>'
    # Escape '$' so Regex.Replace doesn't interpret it as a backreference token
    # in the replacement pattern.
    $newLiteralEscaped = $newLiteral.Replace('$', '$$$$')

    return [System.Text.RegularExpressions.Regex]::Replace($Xml, $pattern, $newLiteralEscaped)
}

$full = Replace-RunText $full `
    "SmartCash là một mô hình quản trị cộng đồng hợp tác và tăng trưởng tập trung vào tiền tệ và một nền kinh tế phi tập trung dựa trên blockchain." `
    "SmartCash được quản trị bởi cộng đồng nhằm tạo nên một phương tiện thanh toán và một nền kinh tế phi tập trung dựa trên blockchain."

$full = Replace-RunText $full `
    "Chúng tôi cố gắng cho phép cộng đồng kiểm soát số phận của đồng xu, quản trị, ngân sách, sử dụng và khuyến khích phát triển cộng đồng." `
    "Chúng tôi cố gắng cho phép cộng đồng kiểm soát số phận của đồng coin, từ quản trị, ngân sách, đến việc sử dụng."

$full = Replace-RunText $full `
    "CỘNG ĐỒNG THÚC ĐẨY TĂNG TRƯỞNG" `
    "CỘNG ĐỒNG THÚC ĐẨY"

$full = Replace-RunText $full `
    "TÀI TRỢ" `
    "TÀI TRỢ TĂNG TRƯỞNG"

$full = Replace-RunText $full `
    "Với SmartCash, chúng tôi đã tập trung rất nhiều vào cộng đồng, để dành 70% ngân sách cho các dự án mà cộng đồng muốn theo đuổi, chỉ còn lại 30% cho một số dự án như bảo trì hệ thống và một số nhu cầu cơ bản khác. Chúng tôi đặc biệt khuyến khích bất kỳ ai có kỹ năng giúp SmartCash phát triển để tham gia vào cộng đồng và mang lại những ý tưởng sáng tạo tốt nhất cho bạn." `
    "Với SmartCash, chúng tôi đã tập trung rất nhiều vào cộng đồng, để dành 70% ngân sách cho các dự án mà cộng đồng muốn theo đuổi, chỉ còn lại 30% cho một số dự án như bảo trì hệ thống và một số nhu cầu cơ bản khác. Chúng tôi đặc biệt khuyến khích bất kỳ ai có kỹ năng, ý tưởng sáng tạo giúp SmartCash phát triển để tham gia vào cộng đồng."

$full = Replace-RunText $full `
    "Khai thác mỏ SmartCash ngăn tập trung khai thác mỏ và kích thích phát triển mạng lưới. Mỗi máy tính có thể được sử dụng như một thiết bị khai thác mỏ trong khi vẫn cho phép máy tính đó được sử dụng cho các nhiệm vụ khác. ASICs have yet to be created for the Keccak mining algorithm and it’s probably safe to assume no ASICs will be created for quite some time." `
    "Khai thác SmartCash ngăn chặn khai thác tập trung và kích thích phát triển mạng lưới. Mỗi máy tính có thể được sử dụng như một thiết bị khai thác mỏ trong khi vẫn cho phép máy tính đó được sử dụng cho các nhiệm vụ khác. ASICs have yet to be created for the Keccak mining algorithm and it’s probably safe to assume no ASICs will be created for quite some time."

$full = Replace-RunText $full `
    "Toàn bộ nền tảng của SmartCash dựa trên sự tăng trưởng việc chấp nhận của cộng đồng. Mọi người đều là thành viên của nhóm SmartCash." `
    "Toàn bộ nền tảng của SmartCash dựa trên việc chấp nhận của cộng đồng. Mọi người đều là thành viên của nhóm SmartCash."

$d.Content.InsertXML($full)
